$wb = $excel.ActiveWorkbook

# Hunk 0: ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 6186.467
$ws.Range("J112").Value = 1621.2142
$ws.Range("L112").Value = 4863.642599999999
$ws.Range("N112").Value = -7079.642599999999

# Hunk 1: ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1036.5
$ws.Range("I132").Value = 971.36365
$ws.Range("J132").Value = 1753
$ws.Range("K132").Value = 2914.09095
$ws.Range("L132").Value = 5259
$ws.Range("M132").Value = -384.0909499999998
$ws.Range("N132").Value = -10319

# Hunk 2: ALC row 133
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 55857.145
$ws.Range("J133").Value = 55857.145
$ws.Range("L133").Value = 55857.145
$ws.Range("N133").Value = -65977.14499999999

# Hunk 3: ALC row 134
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 67945
$ws.Range("J134").Value = 67945
$ws.Range("L134").Value = 67945
$ws.Range("N134").Value = -78085

# Hunk 4: ALC row 136
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 70803.5
$ws.Range("J136").Value = 70803.5
$ws.Range("L136").Value = 70803.5
$ws.Range("N136").Value = -81003.5

# Hunk 5: ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2800.5789
$ws.Range("I138").Value = 1410.4736
$ws.Range("J138").Value = 4190.684
$ws.Range("K138").Value = 4231.4208
$ws.Range("L138").Value = 12572.052
$ws.Range("M138").Value = 908.5792000000001
$ws.Range("N138").Value = -22852.052

# Hunk 6: ALC row 140
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 245950
$ws.Range("J140").Value = 245950
$ws.Range("L140").Value = 245950
$ws.Range("N140").Value = -256310

# Hunk 7: ARM row 15
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 22102.5
$ws.Range("I15").Value = 1000
$ws.Range("J15").Value = 24447.223
$ws.Range("K15").Value = 1000
$ws.Range("L15").Value = 24447.223
$ws.Range("M15").Value = -650
$ws.Range("N15").Value = -25147.223

# Hunk 8: ARM row 31
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 16490.166
$ws.Range("I31").Value = 16490.166
$ws.Range("K31").Value = 16490.166
$ws.Range("M31").Value = -16196.166

# Hunk 9: ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3755.9524
$ws.Range("I132").Value = 3761.8572
$ws.Range("J132").Value = 3744.1428
$ws.Range("K132").Value = 11285.5716
$ws.Range("L132").Value = 11232.4284
$ws.Range("M132").Value = -8755.571599999999
$ws.Range("N132").Value = -16292.4284

# Hunk 10: BSM row 132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 50615.383
$ws.Range("J132").Value = 50615.383
$ws.Range("L132").Value = 50615.383
$ws.Range("N132").Value = -60735.383

# Hunk 11: CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 503
$ws.Range("I107").Value = 341.46667
$ws.Range("J107").Value = 805.875
$ws.Range("K107").Value = 341.46667
$ws.Range("L107").Value = 805.875
$ws.Range("M107").Value = 1578.53333
$ws.Range("N107").Value = -4645.875

# Hunk 12: CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4812.923
$ws.Range("I122").Value = 2110.5715
$ws.Range("J122").Value = 11691.637
$ws.Range("K122").Value = 6331.7145
$ws.Range("L122").Value = 35074.911
$ws.Range("M122").Value = -3881.7145
$ws.Range("N122").Value = -39974.911

# Hunk 13: CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3828.375
$ws.Range("I132").Value = 4023.2273
$ws.Range("J132").Value = 3399.7
$ws.Range("K132").Value = 12069.6819
$ws.Range("L132").Value = 10199.1
$ws.Range("M132").Value = -9539.6819
$ws.Range("N132").Value = -15259.1

# Hunk 14: CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2331.8647
$ws.Range("I134").Value = 1946.5862
$ws.Range("K134").Value = 5839.7586
$ws.Range("M134").Value = -3304.7586

# Hunk 15: CUL row 9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 69521.42999999999
$ws.Range("J9").Value = 81000
$ws.Range("L9").Value = 243000
$ws.Range("N9").Value = -243448

# Hunk 16: CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2410.913
$ws.Range("I68").Value = 863.5135
$ws.Range("J68").Value = 4200.0938
$ws.Range("K68").Value = 2590.5405
$ws.Range("L68").Value = 12600.2814
$ws.Range("M68").Value = -1779.5405
$ws.Range("N68").Value = -14222.2814

# Hunk 17: CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2410.913
$ws.Range("I71").Value = 863.5135
$ws.Range("J71").Value = 4200.0938
$ws.Range("K71").Value = 7771.6215
$ws.Range("L71").Value = 37800.8442
$ws.Range("M71").Value = -3715.6215
$ws.Range("N71").Value = -45912.8442

# Hunk 18: CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 741.5088
$ws.Range("J107").Value = 1888
$ws.Range("L107").Value = 5664
$ws.Range("N107").Value = -9504

# Hunk 19: CUL row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 1995
$ws.Range("I136").Value = 851.1
$ws.Range("J136").Value = 2709.9375
$ws.Range("K136").Value = 2553.3
$ws.Range("L136").Value = 8129.8125
$ws.Range("M136").Value = 2546.7
$ws.Range("N136").Value = -18329.8125

# Hunk 20: GSM row 14
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 107100136
$ws.Range("I14").Value = 107100136
$ws.Range("K14").Value = 107100136
$ws.Range("M14").Value = -107099968

# Hunk 21: GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1383.3334
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1383.3334
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").Value = 4150.0002
$ws.Range("N122").Value = -9050.0002

# Hunk 22: GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6226.2
$ws.Range("I132").Value = 2584
$ws.Range("J132").Value = 9868.4
$ws.Range("K132").Value = 7752
$ws.Range("L132").Value = 29605.2
$ws.Range("M132").Value = -5222
$ws.Range("N132").Value = -34665.2

# Hunk 23: LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2296.125
$ws.Range("I7").Value = 2296.125
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2296.125
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -2184.125

# Hunk 24: LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2629.6428
$ws.Range("I40").Value = 2523.88
$ws.Range("J40").Value = 3511
$ws.Range("K40").Value = 2523.88
$ws.Range("L40").Value = 3511
$ws.Range("M40").Value = -2387.88
$ws.Range("N40").Value = -3783

# Hunk 25: LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3116.6667
$ws.Range("I68").Value = 2550
$ws.Range("J68").Value = 4250
$ws.Range("K68").Value = 2550
$ws.Range("L68").Value = 4250
$ws.Range("M68").Value = -1801
$ws.Range("N68").Value = -5748

# Hunk 26: LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3116.6667
$ws.Range("I71").Value = 2550
$ws.Range("J71").Value = 4250
$ws.Range("K71").Value = 12750
$ws.Range("L71").Value = 21250
$ws.Range("M71").Value = -9006
$ws.Range("N71").Value = -28738

# Hunk 27: LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5369.769
$ws.Range("I122").Value = 5528.8125
$ws.Range("J122").Value = 4642.7144
$ws.Range("K122").Value = 16586.4375
$ws.Range("L122").Value = 13928.1432
$ws.Range("M122").Value = -14136.4375
$ws.Range("N122").Value = -18828.1432

# Hunk 28: LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2296.125
$ws.Range("I126").Value = 2296.125
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6888.375
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -4418.375

# Hunk 29: LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2743.634
$ws.Range("I132").Value = 1799.9048
$ws.Range("J132").Value = 3734.55
$ws.Range("K132").Value = 5399.7144
$ws.Range("L132").Value = 11203.65
$ws.Range("M132").Value = -2869.7144
$ws.Range("N132").Value = -16263.65

# Hunk 30: WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1585.5714
$ws.Range("I122").Value = 1649.8334
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 4949.5002
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -2499.5002
$ws.Range("N122").Value = -8500

# Hunk 31: WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1513.3334
$ws.Range("I126").Value = 1513.3334
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4540.0002
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -2070.0002

# Hunk 32: WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4315.7646
$ws.Range("I132").Value = 4796.1113
$ws.Range("J132").Value = 3775.375
$ws.Range("K132").Value = 14388.3339
$ws.Range("L132").Value = 11326.125
$ws.Range("M132").Value = -11858.3339
$ws.Range("N132").Value = -16386.125

# Hunk 33: WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4744.66
$ws.Range("I136").Value = 1510.8214
$ws.Range("J136").Value = 8860.454
$ws.Range("K136").Value = 4532.4642
$ws.Range("L136").Value = 26581.362
$ws.Range("M136").Value = -1982.4642
$ws.Range("N136").Value = -31681.362

# Hunk 34: WVR row 138
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 38701.285
$ws.Range("J138").Value = 38701.285
$ws.Range("L138").Value = 38701.285
$ws.Range("N138").Value = -48981.285
